# Auto-generated edit script for resum_diari_meteocat.xlsx
# Applies the 2026-02-27 19:50 automatic data/banner update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-27 19:48:43"
$ws.Range("H2").Value = "'53%"
$ws.Range("E3").Value = "2026-02-27 19:48:46"
$ws.Range("E4").Value = "2026-02-27 19:48:48"
$ws.Range("E5").Value = "2026-02-27 19:48:51"
$ws.Range("H5").Value = "'40%"
$ws.Range("I5").Value = "0.0 mm"
$ws.Range("N5").Value = "1.4 °C 19:29 TU"
$ws.Range("E6").Value = "2026-02-27 19:48:54"
$ws.Range("E7").Value = "2026-02-27 19:48:57"
$ws.Range("H7").Value = "'86%"
$ws.Range("J7").Value = "1024.8 hPa"
$ws.Range("E8").Value = "2026-02-27 19:49:00"
$ws.Range("H8").Value = "'63%"
$ws.Range("N8").Value = "8.4 °C 19:27 TU"
$ws.Range("O8").Value = "12.0 °C"
$ws.Range("E9").Value = "2026-02-27 19:49:02"
$ws.Range("H9").Value = "'90%"
$ws.Range("E10").Value = "2026-02-27 19:49:05"
$ws.Range("E11").Value = "2026-02-27 19:49:08"
$ws.Range("E12").Value = "2026-02-27 19:49:10"
$ws.Range("E13").Value = "2026-02-27 19:49:13"
$ws.Range("G13").Value = "3 cm"
$ws.Range("J13").Value = "1025.5 hPa"
$ws.Range("E14").Value = "2026-02-27 19:49:16"
$ws.Range("E15").Value = "2026-02-27 19:49:18"
$ws.Range("E16").Value = "2026-02-27 19:49:21"
$ws.Range("O16").Value = "2.8 °C"
$ws.Range("E17").Value = "2026-02-27 19:49:24"
$ws.Range("N17").Value = "5.1 °C 19:06 TU"
$ws.Range("O17").Value = "7.7 °C"
$ws.Range("E18").Value = "2026-02-27 19:49:26"
$ws.Range("E19").Value = "2026-02-27 19:49:29"
$ws.Range("O19").Value = "10.5 °C"
$ws.Range("E20").Value = "2026-02-27 19:49:31"
$ws.Range("E21").Value = "2026-02-27 19:49:34"
$ws.Range("H21").Value = "'59%"
$ws.Range("J21").Value = "1024.3 hPa"
$ws.Range("O21").Value = "9.9 °C"
$ws.Range("E22").Value = "2026-02-27 19:49:37"
$ws.Range("O22").Value = "1.5 °C"
$ws.Range("E23").Value = "2026-02-27 19:49:40"
$ws.Range("N23").Value = "1.4 °C 19:06 TU"
$ws.Range("O23").Value = "3.8 °C"
$ws.Range("E24").Value = "2026-02-27 19:49:42"
$ws.Range("J24").Value = "1023.5 hPa"
$ws.Range("E25").Value = "2026-02-27 19:49:45"
$ws.Range("E26").Value = "2026-02-27 19:49:48"
$ws.Range("O26").Value = "10.5 °C"
$ws.Range("E27").Value = "2026-02-27 19:49:50"
$ws.Range("N27").Value = "3.2 °C 19:29 TU"
$ws.Range("E28").Value = "2026-02-27 19:49:53"
$ws.Range("E29").Value = "2026-02-27 19:49:56"
$ws.Range("E30").Value = "2026-02-27 19:49:59"
$ws.Range("E31").Value = "2026-02-27 19:50:01"
$ws.Range("E32").Value = "2026-02-27 19:50:04"
$ws.Range("H32").Value = "'58%"
$ws.Range("E33").Value = "2026-02-27 19:50:07"
$ws.Range("H33").Value = "'51%"
$ws.Range("J33").Value = "1023.7 hPa"
$ws.Range("E34").Value = "2026-02-27 19:50:10"
$ws.Range("H34").Value = "'47%"
$ws.Range("E35").Value = "2026-02-27 19:50:13"
$ws.Range("O35").Value = "12.1 °C"
$ws.Range("E36").Value = "2026-02-27 19:50:15"
$ws.Range("J36").Value = "1024.8 hPa"
$ws.Range("E37").Value = "2026-02-27 19:50:18"
$ws.Range("H37").Value = "'68%"
$ws.Range("E38").Value = "2026-02-27 19:50:21"
$ws.Range("O38").Value = "10.3 °C"
$ws.Range("E39").Value = "2026-02-27 19:50:24"
$ws.Range("I39").Value = "0.0 mm"
$ws.Range("E40").Value = "2026-02-27 19:50:27"
$ws.Range("J40").Value = "1024.8 hPa"
$ws.Range("E41").Value = "2026-02-27 19:50:29"
$ws.Range("E42").Value = "2026-02-27 19:50:32"
$ws.Range("O42").Value = "11.5 °C"
$ws.Range("E43").Value = "2026-02-27 19:50:34"
$ws.Range("E44").Value = "2026-02-27 19:50:37"
$ws.Range("E45").Value = "2026-02-27 19:50:40"
$ws.Range("H45").Value = "'42%"
$ws.Range("E46").Value = "2026-02-27 19:50:42"
